$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.432.66"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "2.288.14"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "156.86"
$ws.Range("E5").Value = "  +15,570.89%  "
$ws.Range("D6").Value = "306.86"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").Value = "95.99"
$ws.Range("E7").Value = "  +4.61%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "0.496"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").Value = "36.01"
$ws.Range("E11").Value = "  +11.11%  "
$ws.Range("D12").Value = "0.0805"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "2.641.83"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "14.51"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "2.307.06"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").Value = "0.802"
$ws.Range("E18").Value = "  +5.53%  "
$ws.Range("D19").Value = "42.331.84"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("D21").Value = "0.0₃0920"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").Value = "6.01"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").Value = "68.17"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").Value = "243.16"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "2.61"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").Value = "1.95"
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "24.13"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "36.15"
$ws.Range("E29").Value = "  +6.34%  "
$ws.Range("D30").Value = "9.64"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").Value = "2.10"
$ws.Range("E31").Value = "  -8.71%  "
$ws.Range("D32").Value = "161.00"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "5.35"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").Value = "3.08"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("E37").Value = "  +4.82%  "
$ws.Range("D38").Value = "17.19"
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "4.20"
$ws.Range("E42").Value = "  +7.83%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "19.69"
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.013.59"
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("E45").Value = "  +10.85%  "
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "3.00"
$ws.Range("E47").Value = "  +4.26%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "10.19"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").Value = "53.44"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("D51").Value = "73.13"
$ws.Range("E51").Value = "  +0.21%  "
